# Update rich cards for button template
# - Adds a new PROMPTS_CONFIG row (ESI_PHA_WELCOME_MSG / welcome text / TEXT)
# - Adds two new RICH_CARDS_CONFIG rows (BUTTONS: WISMO + Payment/Outstanding Balance)
# - Updates the view state (selection / active sheet / zoom) on both sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PROMPTS_CONFIG")
$ws2 = $wb.Worksheets.Item("RICH_CARDS_CONFIG")

# NOTE: cell values are written in the precise order needed so that newly
# interned shared-strings land at the same table indices as the target file
# (new unique strings are appended to xl/sharedStrings.xml in write order).

# --- PROMPTS_CONFIG: new row 19 (RESPONSE_ID / MEDIA_TYPE only first) ----
$ws1.Range("A19").Value = "ESI_PHA_WELCOME_MSG"
$ws1.Range("E19").Value = "TEXT"

# --- RICH_CARDS_CONFIG: new rows 7 & 8 -----------------------------------
$ws2.Range("A7").Value = "ESI_PHA_WELCOME_MSG"
$ws2.Range("B7").Value = "BUTTONS"
$ws2.Range("C7").Value = "WEB"
$ws2.Range("D7").Value = "ESA_PHA_WISMO"
$ws2.Range("E7").Value = "Where is my order?"

$ws2.Range("A8").Value = "ESI_PHA_WELCOME_MSG"
$ws2.Range("B8").Value = "BUTTONS"
$ws2.Range("C8").Value = "WEB"
$ws2.Range("D8").Value = "ESA_PHA_PAYMENT_BAL"
$ws2.Range("E8").Value = "Payment/Outstanding Balance"

# --- PROMPTS_CONFIG: welcome message text (new string, written last) ----
$ws1.Range("B19").Value = "Hello, I’m Ask Iris, Express Scripts digital assistant. I'm in beta mode and still in learning phase. Here’s how I can help you right now:"

# --- View state: selection, active sheet, zoom ---------------------------
$ws1.Range("C19").Select()

$ws2.Activate()
$ws2.Range("E8").Select()
$excel.ActiveWindow.Zoom = 64
